$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new tracked video row (row 4) with its translation-time stats,
# mirroring the layout/styles used by the existing rows.
$ws.Range("A4").Value = "makeine vol2.1"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 45
$ws.Range("E4").Value = 331

$ws.Range("F4").Value = 4569
$ws.Range("F4").NumberFormat = $ws.Range("F2").NumberFormat

$ws.Range("G4").Formula = "=F4/E4"
$ws.Range("G4").NumberFormat = $ws.Range("G2").NumberFormat

$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 0.92

$ws.Range("J4").Formula = "=I4/F4"
$ws.Range("J4").NumberFormat = $ws.Range("J2").NumberFormat

# Match the author's final selection/view state.
$ws.Range("J4").Select() | Out-Null
